$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - rows 2-11 map to F column values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 8092
$wsExhibit.Range("F3").Value = 82
$wsExhibit.Range("F4").Value = 225
$wsExhibit.Range("F5").Value = 67
$wsExhibit.Range("F6").Value = 833
$wsExhibit.Range("F7").Value = 1352
$wsExhibit.Range("F8").Value = 221
$wsExhibit.Range("F9").Value = 27
$wsExhibit.Range("F10").Value = 196
$wsExhibit.Range("F11").Value = 53

# Sheet "全部类型" (all types) - rows 2-12 map to F column values (row 9 unchanged)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 8092
$wsAll.Range("F3").Value = 82
$wsAll.Range("F4").Value = 225
$wsAll.Range("F5").Value = 67
$wsAll.Range("F6").Value = 833
$wsAll.Range("F7").Value = 1352
$wsAll.Range("F8").Value = 221
$wsAll.Range("F10").Value = 27
$wsAll.Range("F11").Value = 196
$wsAll.Range("F12").Value = 53
